$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("final allocation")

$ws.Range("B2").Value = 223
$ws.Range("B3").Value = 202
$ws.Range("B4").Value = 222
$ws.Range("B5").Value = 102
$ws.Range("B6").Value = 126
